# Fill in the bill-form header fields (name, designation, year, term,
# department/branch, subject department, amount in words) and resize
# column A / row 36 to fit the newly entered text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header block (rows 3-5) ---------------------------------------------
# A3 label "নাম:" -> append teacher's name
$ws.Range("A3").Value = "নাম: Dr. Md. Alhaz Uddin "

# A4 label "পদবী: " -> append designation
$ws.Range("A4").Value = "পদবী: অধ্যাপক"

# G4 / I4 were blank -> fill academic year / term
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"

# B5 was blank -> branch
$ws.Range("B5").Value = "সিএসই"

# F5 label "বিভাগ :" -> append department
$ws.Range("F5").Value = "বিভাগ :গণিত"

# --- Amount in words (row 32) ---------------------------------------------
$ws.Range("A32").Value = "কথায়:দুই হাজার সাতশো টাকা মাত্র।"

# --- Column / row sizing ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Rows.Item(36).RowHeight = 68.4

# --- View state (best effort; cosmetic only) --------------------------------
$ws.Range("I32").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
